$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 599.2325665869322
$ws.Range("C2").Value = 1198.4651331738753
$ws.Range("D2").Value = 0.000000000005535606695250550047
$ws.Range("E2").Value = 1199.2925665866342
$ws.Range("B3").Value = 550.5681603505575
$ws.Range("C3").Value = 47234.44036835699
$ws.Range("D3").Value = 46133.3040476568
$ws.Range("E3").Value = 1750.6881603497955
$ws.Range("B4").Value = 533.1576321477303
$ws.Range("C4").Value = 66122.97364050394
$ws.Range("D4").Value = 65056.658376208914
$ws.Range("E4").Value = 2333.3376321464652
$ws.Range("B5").Value = 524.6474060310906
$ws.Range("C5").Value = 85977.94637022456
$ws.Range("D5").Value = 84928.65155816046
$ws.Range("E5").Value = 2924.8874060292997
$ws.Range("B6").Value = 520.0998994174328
$ws.Range("C6").Value = 128452.75300706379
$ws.Range("D6").Value = 127412.55320822689
$ws.Range("E6").Value = 3520.3998994151225
$ws.Range("B7").Value = 516.669980906794
$ws.Range("C7").Value = 139827.84571725075
$ws.Range("D7").Value = 138794.50575543573
$ws.Range("E7").Value = 4117.0299809040835
$ws.Range("B8").Value = 514.1937734274638
$ws.Range("C8").Value = 241288.17562609198
$ws.Range("D8").Value = 240259.78807922665
$ws.Range("E8").Value = 4714.613773428269
$ws.Range("B9").Value = 512.4251804644493
$ws.Range("C9").Value = 474203.46256445814
$ws.Range("D9").Value = 473178.612203523
$ws.Range("E9").Value = 5312.9051804688725
$ws.Range("B10").Value = 511.11950089852843
$ws.Range("C10").Value = 278749.9205303257
$ws.Range("D10").Value = 277727.6815285184
$ws.Range("E10").Value = 5911.659500906605
$ws.Range("B11").Value = 509.8739625918634
$ws.Range("C11").Value = 233486.24100547025
$ws.Range("D11").Value = 232466.49308028346
$ws.Range("E11").Value = 6510.47396260371
$ws.Range("B12").Value = 509.08157488863094
$ws.Range("C12").Value = 312684.3155370875
$ws.Range("D12").Value = 311666.1523873043
$ws.Range("E12").Value = 7109.7415749042675
$ws.Range("B13").Value = 508.3617511673978
$ws.Range("C13").Value = 299997.5700394721
$ws.Range("D13").Value = 298980.846537129
$ws.Range("E13").Value = 7709.08175118687
$ws.Range("B14").Value = 507.7011314114033
$ws.Range("C14").Value = 413182.4950619857
$ws.Range("D14").Value = 412167.0927991502
$ws.Range("E14").Value = 8308.481131433082
$ws.Range("B15").Value = 507.1134782262492
$ws.Range("C15").Value = 480509.72463630274
$ws.Range("D15").Value = 479495.49767984013
$ws.Range("E15").Value = 8907.95347824321
$ws.Range("B16").Value = 506.70891039184653
$ws.Range("C16").Value = 458536.8454511972
$ws.Range("D16").Value = 457523.4276304002
$ws.Range("E16").Value = 9507.60891040408
$ws.Range("B17").Value = 506.37419568797947
$ws.Range("C17").Value = 346096.3856426927
$ws.Range("D17").Value = 345083.63725130016
$ws.Range("E17").Value = 10107.334195695348
$ws.Range("B18").Value = 505.9156994124931
$ws.Range("C18").Value = 345101.1809678955
$ws.Range("D18").Value = 344089.34956905077
$ws.Range("E18").Value = 10706.935699415151
$ws.Range("B19").Value = 505.5574485597781
$ws.Range("C19").Value = 274185.05163111107
$ws.Range("D19").Value = 273173.9367339911
$ws.Range("E19").Value = 11306.637448557525
$ws.Range("B20").Value = 505.3848325533055
$ws.Range("C20").Value = 676472.6870760949
$ws.Range("D20").Value = 675461.9174109857
$ws.Range("E20").Value = 11906.524832546196
$ws.Range("B21").Value = 505.02275530307867
$ws.Range("C21").Value = 393241.1453356835
$ws.Range("D21").Value = 392231.0998250692
$ws.Range("E21").Value = 12506.22275529106
$ws.Range("B22").Value = 504.80390114795733
$ws.Range("C22").Value = 396514.2369854308
$ws.Range("D22").Value = 395504.62918312603
$ws.Range("E22").Value = 13106.063901131056
$ws.Range("B23").Value = 504.5788413421824
$ws.Range("C23").Value = 398656.731408491
$ws.Range("D23").Value = 397647.57372579694
$ws.Range("E23").Value = 13705.898841320297
$ws.Range("B24").Value = 504.3789067816731
$ws.Range("C24").Value = 637264.6369126254
$ws.Range("D24").Value = 636255.8790990569
$ws.Range("E24").Value = 14305.758906754849
$ws.Range("B25").Value = 504.21132256873136
$ws.Range("C25").Value = 533212.3259117552
$ws.Range("D25").Value = 532203.9032666121
$ws.Range("E25").Value = 14905.6513225369
$ws.Range("B26").Value = 503.9670816506262
$ws.Range("C26").Value = 490949.3519545602
$ws.Range("D26").Value = 489941.4177912542
$ws.Range("E26").Value = 15505.467081613857
$ws.Range("B27").Value = 503.90830144229466
$ws.Range("C27").Value = 558070.5745679992
$ws.Range("D27").Value = 557062.7579651166
$ws.Range("E27").Value = 16105.468301400564
$ws.Range("B28").Value = 503.64623963061166
$ws.Range("C28").Value = 688757.4085990215
$ws.Range("D28").Value = 687750.1161197624
$ws.Range("E28").Value = 16705.266239593333
$ws.Range("B29").Value = 503.56085757351946
$ws.Range("C29").Value = 637614.4997023202
$ws.Range("D29").Value = 636607.3779871752
$ws.Range("E29").Value = 17305.24085754887
$ws.Range("B30").Value = 503.433130893882
$ws.Range("C30").Value = 551522.3639572252
$ws.Range("D30").Value = 550515.497695438
$ws.Range("E30").Value = 17905.1731308819
$ws.Range("B31").Value = 503.3760939687362
$ws.Range("C31").Value = 985098.8282358274
$ws.Range("D31").Value = 984092.0760478931
$ws.Range("E31").Value = 18505.17609396946